$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:J102")
$key = $ws.Range("A2")
$sortRange.Sort($key, 1, [Type]::Missing, [Type]::Missing, 1, [Type]::Missing, 1, 1)

$rng = $ws.Range("A2:J102")
$arr = @("1.8V GPIO")
$rng.AutoFilter(9, $arr, 7)
Write-Host "done"
